$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 197, shifting existing rows 197:286 down to 198:287
$ws.Rows.Item(197).Insert()

# Populate the newly inserted row 197 with the new data record
$ws.Cells.Item(197, 1).Value = 4
$ws.Cells.Item(197, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(197, 3).Value = "Los Lagos"
$ws.Cells.Item(197, 4).Value = 44636
$ws.Cells.Item(197, 4).NumberFormat = $ws.Cells.Item(198, 4).NumberFormat
$ws.Cells.Item(197, 5).Value = 10
$ws.Cells.Item(197, 6).Value = 100112008
$ws.Cells.Item(197, 7).Value = "Coliflor"
$ws.Cells.Item(197, 8).Value = "Sin especificar"
$ws.Cells.Item(197, 9).Value = "Primera"
$ws.Cells.Item(197, 10).Value = 100
$ws.Cells.Item(197, 11).Value = 1600
$ws.Cells.Item(197, 12).Value = 1600
$ws.Cells.Item(197, 13).Value = 1600
$ws.Cells.Item(197, 14).Value = "$/unidad"
$ws.Cells.Item(197, 15).Value = "Región Metropolitana"
$ws.Cells.Item(197, 16).Value = 1600
$ws.Cells.Item(197, 17).Value = 1
$ws.Cells.Item(197, 18).Value = "Hortaliza"
